# Weekly update: a new price-survey row for "Ají" (Feria Lagunitas de Puerto
# Montt) is inserted just before the existing row 208, pushing every
# subsequent data row (old 208..314) down by one (new 209..315).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new, blank row at position 208 (entire row), shifting rows 208-314
# down to 209-315 and extending the used range to A1:R315.
$ws.Range("A208:R208").EntireRow.Insert()

# Populate the newly inserted row 208 with the new survey data.
$ws.Cells.Item(208, 1).Value  = 4
$ws.Cells.Item(208, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(208, 3).Value  = "Los Lagos"
$ws.Cells.Item(208, 4).Value  = 44830
$ws.Cells.Item(208, 5).Value  = 10
$ws.Cells.Item(208, 6).Value  = 100112021
$ws.Cells.Item(208, 7).Value  = "Ají"
$ws.Cells.Item(208, 8).Value  = "Inferno"
$ws.Cells.Item(208, 9).Value  = "Primera"
$ws.Cells.Item(208, 10).Value = 70
$ws.Cells.Item(208, 11).Value = 23000
$ws.Cells.Item(208, 12).Value = 23000
$ws.Cells.Item(208, 13).Value = 23000
$ws.Cells.Item(208, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(208, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(208, 16).Value = 2300
$ws.Cells.Item(208, 17).Value = 10
$ws.Cells.Item(208, 18).Value = "Hortaliza"
